$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card10")

# Row 21: previously-blank tracking columns (B-K, M) get the literal
# placeholder text "nan", matching the pattern used by every other
# event row (14-20) on this sheet. L/N/O already held data and are
# left untouched.
$ws.Range("B21").Value = "nan"
$ws.Range("C21").Value = "nan"
$ws.Range("D21").Value = "nan"
$ws.Range("E21").Value = "nan"
$ws.Range("F21").Value = "nan"
$ws.Range("G21").Value = "nan"
$ws.Range("H21").Value = "nan"
$ws.Range("I21").Value = "nan"
$ws.Range("J21").Value = "nan"
$ws.Range("K21").Value = "nan"
$ws.Range("M21").Value = "nan"

# Row 22: new service event appended to the bottom of the Card10 log.
# A22 mirrors the "card" number stored as text (like every other cell in
# column A on this sheet), so force Text format before assigning the
# numeric-looking string to stop Excel auto-typing it as a number.
$ws.Range("A22").NumberFormat = "@"
$ws.Range("A22").Value = "10"
$ws.Range("L22").Value = '6\3\2025'
$ws.Range("N22").Value = "تم سن الفلاتس وتغيير الجرائد الخلفيه (1_5_8)"
$ws.Range("O22").Value = "الخبير"
